# Weekly data refresh: insert this week's new record at the top of the
# Agrícola del Norte S.A. de Arica - Espinaca data block (row 80), pushing
# the rest of the history down by one row (oldest row falls to row 114).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

$ws.Range("A80").Value = 1
$ws.Range("B80").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C80").Value = 'Arica y Parinacota'
$ws.Range("D80").Value = 45202
$ws.Range("E80").Value = 15
$ws.Range("F80").Value = 100112012
$ws.Range("G80").Value = 'Espinaca'
$ws.Range("H80").Value = 'Sin especificar'
$ws.Range("I80").Value = 'Primera'
$ws.Range("J80").Value = 270
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = 1900
$ws.Range("N80").Value = '$/atado 2,5 a 3 kilos'
$ws.Range("O80").Value = 'Región de Arica y Parinacota'
$ws.Range("P80").Value = 633
$ws.Range("Q80").Value = 3
$ws.Range("R80").Value = 'Hortaliza'
